# Atualização de bases das ligas, do dia: 21-04-2024 às 13:33
#
# Applies the source-data refresh to the "Israel Premier League" sheet:
#  - 4 pairs of rows had their (non-id) fields swapped back into the
#    correct id-order (24/25, 44/45, 72/73, 108/110)
#  - rows 201/202 got refreshed odds data
#  - rows 203-207 (stale duplicate fixtures) were removed entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row pair 24 <-> 25 (column A / id stays put) ----
$ws.Range("B24").Value = 6799846
$ws.Range("F24").Value = "Hapoel Jerusalem FC"
$ws.Range("G24").Value = "Maccabi Netanya"
$ws.Range("K24").Value = 2.8
$ws.Range("L24").Value = 3.3
$ws.Range("M24").Value = 2.3
$ws.Range("N24").Value = 2.5
$ws.Range("O24").Value = 3.2
$ws.Range("P24").Value = 2.6
$ws.Range("R24").Value = 1.9
$ws.Range("S24").Value = 1.95
$ws.Range("U24").Value = 2.05
$ws.Range("V24").Value = 1.8
$ws.Range("X24").Value = 2.2
$ws.Range("AC24").Value = 0.8

$ws.Range("B25").Value = 6799841
$ws.Range("F25").Value = "MS Ashdod"
$ws.Range("G25").Value = "Maccabi Petach Tikva"
$ws.Range("K25").Value = 2.25
$ws.Range("L25").Value = 3.25
$ws.Range("M25").Value = 2.75
$ws.Range("N25").Value = 2.3
$ws.Range("O25").Value = 3.3
$ws.Range("P25").Value = 2.7
$ws.Range("R25").Value = 1.75
$ws.Range("S25").Value = 2.05
$ws.Range("U25").Value = 2
$ws.Range("V25").Value = 1.85
$ws.Range("X25").Value = 2.3
$ws.Range("AC25").Value = 0.8500000000000001

# ---- Row pair 44 <-> 45 ----
$ws.Range("B44").Value = 7511180
$ws.Range("F44").Value = "Maccabi Petach Tikva"
$ws.Range("G44").Value = "Maccabi Netanya"
$ws.Range("H44").Value = 1
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = "H"
$ws.Range("K44").Value = 2.5
$ws.Range("L44").Value = 3.2
$ws.Range("M44").Value = 2.5
$ws.Range("N44").Value = 2.5
$ws.Range("O44").Value = 3.2
$ws.Range("P44").Value = 2.45
$ws.Range("R44").Value = 1.95
$ws.Range("S44").Value = 1.9
$ws.Range("T44").Value = 2.5
$ws.Range("U44").Value = 2
$ws.Range("V44").Value = 1.85
$ws.Range("W44").Value = 1.5
$ws.Range("Y44").Value = -1
$ws.Range("Z44").Value = 0.95
$ws.Range("AA44").Value = -1
$ws.Range("AB44").Value = -1
$ws.Range("AC44").Value = 0.8500000000000001

$ws.Range("B45").Value = 7511181
$ws.Range("F45").Value = "Hapoel Hadera"
$ws.Range("G45").Value = "Maccabi Bnei Raina"
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 2
$ws.Range("J45").Value = "A"
$ws.Range("K45").Value = 2.7
$ws.Range("L45").Value = 3.3
$ws.Range("M45").Value = 2.4
$ws.Range("N45").Value = 2.45
$ws.Range("O45").Value = 3.25
$ws.Range("P45").Value = 2.7
$ws.Range("R45").Value = 1.825
$ws.Range("S45").Value = 2.025
$ws.Range("T45").Value = 2.25
$ws.Range("U45").Value = 1.875
$ws.Range("V45").Value = 1.975
$ws.Range("W45").Value = -1
$ws.Range("Y45").Value = 1.7
$ws.Range("Z45").Value = -1
$ws.Range("AA45").Value = 1.025
$ws.Range("AB45").Value = -0.5
$ws.Range("AC45").Value = 0.4875

# ---- Row pair 72 <-> 73 ----
$ws.Range("B72").Value = 7542719
$ws.Range("F72").Value = "Hapoel Haifa"
$ws.Range("G72").Value = "Maccabi Netanya"
$ws.Range("H72").Value = 2
$ws.Range("J72").Value = "H"
$ws.Range("K72").Value = 2.6
$ws.Range("L72").Value = 3.1
$ws.Range("M72").Value = 2.6
$ws.Range("N72").Value = 2.9
$ws.Range("O72").Value = 3.2
$ws.Range("P72").Value = 2.3
$ws.Range("Q72").Value = 0.25
$ws.Range("R72").Value = 1.8
$ws.Range("S72").Value = 2.05
$ws.Range("T72").Value = 2.5
$ws.Range("U72").Value = 2
$ws.Range("V72").Value = 1.85
$ws.Range("W72").Value = 1.9
$ws.Range("Y72").Value = -1
$ws.Range("Z72").Value = 0.8
$ws.Range("AA72").Value = -1
$ws.Range("AB72").Value = 1
$ws.Range("AC72").Value = -1

$ws.Range("B73").Value = 7542640
$ws.Range("F73").Value = "MS Ashdod"
$ws.Range("G73").Value = "Hapoel Bnei Sakhnin"
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = "A"
$ws.Range("K73").Value = 2.05
$ws.Range("L73").Value = 3.2
$ws.Range("M73").Value = 3.5
$ws.Range("N73").Value = 2.15
$ws.Range("O73").Value = 3.1
$ws.Range("P73").Value = 3.2
$ws.Range("Q73").Value = -0.25
$ws.Range("R73").Value = 1.925
$ws.Range("S73").Value = 1.925
$ws.Range("T73").Value = 2.25
$ws.Range("U73").Value = 1.9
$ws.Range("V73").Value = 1.95
$ws.Range("W73").Value = -1
$ws.Range("Y73").Value = 2.2
$ws.Range("Z73").Value = -1
$ws.Range("AA73").Value = 0.925
$ws.Range("AB73").Value = -1
$ws.Range("AC73").Value = 0.95

# ---- Row pair 108 <-> 110 (row 109 is untouched, in between) ----
$ws.Range("B108").Value = 7542736
$ws.Range("F108").Value = "Hapoel Jerusalem FC"
$ws.Range("G108").Value = "Hapoel Bnei Sakhnin"
$ws.Range("I108").Value = 0
$ws.Range("J108").Value = "D"
$ws.Range("K108").Value = 2.2
$ws.Range("L108").Value = 3.4
$ws.Range("M108").Value = 3.2
$ws.Range("N108").Value = 2.375
$ws.Range("O108").Value = 3.1
$ws.Range("P108").Value = 3.1
$ws.Range("Q108").Value = -0.25
$ws.Range("R108").Value = 2.05
$ws.Range("S108").Value = 1.8
$ws.Range("T108").Value = 2
$ws.Range("U108").Value = 1.925
$ws.Range("V108").Value = 1.925
$ws.Range("X108").Value = 2.1
$ws.Range("Y108").Value = -1
$ws.Range("Z108").Value = -0.5
$ws.Range("AA108").Value = 0.4
$ws.Range("AC108").Value = 0.925

$ws.Range("B110").Value = 7542737
$ws.Range("F110").Value = "MS Ashdod"
$ws.Range("G110").Value = "Hapoel Haifa"
$ws.Range("I110").Value = 1
$ws.Range("J110").Value = "A"
$ws.Range("K110").Value = 3
$ws.Range("L110").Value = 3.2
$ws.Range("M110").Value = 2.45
$ws.Range("N110").Value = 3.2
$ws.Range("O110").Value = 3.25
$ws.Range("P110").Value = 2.3
$ws.Range("Q110").Value = 0.25
$ws.Range("R110").Value = 1.85
$ws.Range("S110").Value = 2
$ws.Range("T110").Value = 2.25
$ws.Range("U110").Value = 1.875
$ws.Range("V110").Value = 1.975
$ws.Range("X110").Value = -1
$ws.Range("Y110").Value = 1.3
$ws.Range("Z110").Value = -1
$ws.Range("AA110").Value = 1
$ws.Range("AC110").Value = 0.9750000000000001

# ---- Row 201: refreshed odds for Hapoel Hadera vs MS Ashdod ----
$ws.Range("B201").Value = 8016154
$ws.Range("E201").Value = 45403.54166666666
$ws.Range("F201").Value = "Hapoel Hadera"
$ws.Range("G201").Value = "MS Ashdod"
$ws.Range("K201").Value = 2.45
$ws.Range("L201").Value = 3.2
$ws.Range("M201").Value = 2.875
$ws.Range("N201").Value = 3.2
$ws.Range("O201").Value = 3.1
$ws.Range("P201").Value = 2.25
$ws.Range("Q201").Value = 0.25
$ws.Range("R201").Value = 1.85
$ws.Range("S201").Value = 2
$ws.Range("T201").Value = 2
$ws.Range("U201").Value = 1.8
$ws.Range("V201").Value = 2.05

# ---- Row 202: refreshed odds for Hapoel Petah Tikva vs Beitar Jerusalem ----
$ws.Range("B202").Value = 8016169
$ws.Range("E202").Value = 45403.59375
$ws.Range("F202").Value = "Hapoel Petah Tikva"
$ws.Range("G202").Value = "Beitar Jerusalem"
$ws.Range("K202").Value = 3.5
$ws.Range("L202").Value = 3.4
$ws.Range("M202").Value = 1.95
$ws.Range("N202").Value = 4
$ws.Range("O202").Value = 3.6
$ws.Range("P202").Value = 1.75
$ws.Range("Q202").Value = 0.5
$ws.Range("R202").Value = 2.05
$ws.Range("S202").Value = 1.8
$ws.Range("T202").Value = 2.5
$ws.Range("U202").Value = 1.925
$ws.Range("V202").Value = 1.925

# ---- Remove the 5 stale trailing rows (203-207) ----
$ws.Range("A203:A207").EntireRow.Delete()
